$wb = $excel.ActiveWorkbook

# Deselect current column-range selection residue on Activity sheet and select column H
# (mirrors end-state captured in target workbook: Activity sheet selection becomes H1:H1048576)
$activitySheet = $wb.Worksheets.Item("Activity")
$activitySheet.Columns("H").Select()

# Add the new "Followup" worksheet as the last sheet in the workbook
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$newSheet.Name = "Followup"

# Populate column A (Type) first, then column B (Comments) to match shared-string insertion order
$newSheet.Range("A1").Value = "Type"
$newSheet.Range("A2").Value = "External"
$newSheet.Range("A3").Value = "Internal"

$newSheet.Range("B1").Value = "Comments"
$newSheet.Range("B2").Value = "External Followup"
$newSheet.Range("B3").Value = "Internal Followup"

# Bold the header row
$newSheet.Range("A1:B1").Font.Bold = $true

# Leave selection on the new sheet at C6, matching the final saved view state
$newSheet.Range("C6").Select()
